$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)
$ws.Columns.Item(8).Insert()
$ws.Range("H1").Value = "property_category"
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
Write-Host "done"
